$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells that are being updated so Excel keeps them as text
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "28.741.99"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.804.85"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "231.77"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "0.5946"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("D8").Value = "0.2780"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "0.06858"
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("D10").Value = "23.45"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "0.07552"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "1.810.23"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "0.6287"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "2.050.30"
$ws.Range("D16").Value = "0.000009307"
$ws.Range("E16").Value = "  -7.88%  "
$ws.Range("D17").Value = "75.52"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "28.712.49"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "5.493"
$ws.Range("E19").Value = "  -6.73%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "211.39"
$ws.Range("E21").Value = "  -7.43%  "
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "6.871"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "154.42"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "7.860"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "0.1275"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "16.42"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "1.446"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "0.06208"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "3.786"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "3.759"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").Value = "1.724"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "1.061"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").Value = "0.6437"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "2.726"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").Value = "0.01708"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").Value = "6.431"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").Value = "1.144.17"
$ws.Range("E41").Value = "  -6.09%  "
$ws.Range("D42").Value = "0.8684"
$ws.Range("E42").Value = "  -6.78%  "
$ws.Range("D44").Value = "100.74"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "1.966.51"
$ws.Range("D46").Value = "60.69"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("E47").Value = "  -6.28%  "
$ws.Range("D48").Value = "1.596"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "8.396"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").Value = "0.05464"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "0.4496"
$ws.Range("E51").Value = "  -1.51%  "
